$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85 (shifts existing rows 85..144 down to 86..145)
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new weekly record
$ws.Range("A85").Value = 5
$ws.Range("B85").Value = "Macroferia Regional de Talca"
$ws.Range("C85").Value = "Maule"
$ws.Range("D85").Value = 44596
$ws.Range("D85").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E85").Value = 7
$ws.Range("F85").Value = 100112031
$ws.Range("G85").Value = "Poroto verde"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 100
$ws.Range("K85").Value = 38000
$ws.Range("L85").Value = 38000
$ws.Range("M85").Value = 38000
$ws.Range("N85").Value = "$/saco 25 kilos"
$ws.Range("O85").Value = "Región del Maule"
$ws.Range("P85").Value = 1520
$ws.Range("Q85").Value = 25
$ws.Range("R85").Value = "Hortaliza"
